$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I11 "Sounds" -> "Effects"
$ws.Range("I11").Value = "Effects"

# Remove the old "Effects" remark that used to sit in I12 (fully remove the cell)
$ws.Range("I12").Clear()

# Insert a new row above row 17 (current "Save Data" row) to hold the new
# "Sounds" line item; this shifts "Save Data" and "Total" down by one row
# and automatically adjusts the SUM formulas' ranges.
$ws.Rows("17:17").Insert()

# Carry the formatting from the row above (row 16) into the freshly
# inserted row 17, matching how Excel fills down borders/number formats.
$ws.Range("C16:F16").Copy($ws.Range("C17:F17"))

# Populate the newly inserted row 17 with the Sounds entry
$ws.Range("C17").Value = "Sounds"
$ws.Range("D17").ClearContents()
$ws.Range("E17").Value = 0.5

# Restore the active cell/selection
$ws.Range("I11").Select()
